# Updates Bmp10-Bmpr2 LR-pairs sheet with new TPM-derived values.
# Applies the numeric updates to columns E,F,G,H,I,J,M,N,O,P,Q,R,S,T for rows 2-19
# (columns A-D, K, L are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.178715
$ws.Cells.Item(2, 8).Value = 0.536145
$ws.Cells.Item(2, 9).Value = 0.09904930989061336
$ws.Cells.Item(2, 10).Value = 0.09904930989061336
$ws.Cells.Item(2, 13).Value = 36.81180933333333
$ws.Cells.Item(2, 14).Value = 110.435428
$ws.Cells.Item(2, 15).Value = 0.2598784967371026
$ws.Cells.Item(2, 16).Value = 0.2598784967371026
$ws.Cells.Item(2, 17).Value = 6.578822505006666
$ws.Cells.Item(2, 18).Value = 59.20940254505999
$ws.Cells.Item(2, 19).Value = 0.02574078575722003
$ws.Cells.Item(2, 20).Value = 0.02574078575722003
$ws.Cells.Item(3, 7).Value = 0.178715
$ws.Cells.Item(3, 8).Value = 0.536145
$ws.Cells.Item(3, 9).Value = 0.09904930989061336
$ws.Cells.Item(3, 10).Value = 0.09904930989061336
$ws.Cells.Item(3, 15).Value = 0.1970278712683331
$ws.Cells.Item(3, 16).Value = 0.197027871268333
$ws.Cells.Item(3, 17).Value = 4.987759317866667
$ws.Cells.Item(3, 18).Value = 44.8898338608
$ws.Cells.Item(3, 19).Value = 0.019515474678345
$ws.Cells.Item(3, 20).Value = 0.019515474678345
$ws.Cells.Item(4, 7).Value = 0.178715
$ws.Cells.Item(4, 8).Value = 0.536145
$ws.Cells.Item(4, 9).Value = 0.09904930989061336
$ws.Cells.Item(4, 10).Value = 0.09904930989061336
$ws.Cells.Item(4, 13).Value = 21.95609833333333
$ws.Cells.Item(4, 14).Value = 65.868295
$ws.Cells.Item(4, 15).Value = 0.1550023737603119
$ws.Cells.Item(4, 16).Value = 0.1550023737603119
$ws.Cells.Item(4, 17).Value = 3.923884113641666
$ws.Cells.Item(4, 18).Value = 35.314957022775
$ws.Cells.Item(4, 19).Value = 0.01535287815236581
$ws.Cells.Item(4, 20).Value = 0.01535287815236581
$ws.Cells.Item(5, 7).Value = 0.178715
$ws.Cells.Item(5, 8).Value = 0.536145
$ws.Cells.Item(5, 9).Value = 0.09904930989061336
$ws.Cells.Item(5, 10).Value = 0.09904930989061336
$ws.Cells.Item(5, 13).Value = 13.23098133333333
$ws.Cells.Item(5, 14).Value = 39.692944
$ws.Cells.Item(5, 15).Value = 0.09340609987756826
$ws.Cells.Item(5, 16).Value = 0.09340609987756825
$ws.Cells.Item(5, 17).Value = 2.364574828986666
$ws.Cells.Item(5, 18).Value = 21.28117346088
$ws.Cells.Item(5, 19).Value = 0.009251809732446841
$ws.Cells.Item(5, 20).Value = 0.009251809732446841
$ws.Cells.Item(6, 7).Value = 0.178715
$ws.Cells.Item(6, 8).Value = 0.536145
$ws.Cells.Item(6, 9).Value = 0.09904930989061336
$ws.Cells.Item(6, 10).Value = 0.09904930989061336
$ws.Cells.Item(6, 13).Value = 22.080681
$ws.Cells.Item(6, 14).Value = 66.242043
$ws.Cells.Item(6, 15).Value = 0.1558818838066577
$ws.Cells.Item(6, 16).Value = 0.1558818838066577
$ws.Cells.Item(6, 17).Value = 3.946148904914999
$ws.Cells.Item(6, 18).Value = 35.515340144235
$ws.Cells.Item(6, 19).Value = 0.01543999301549823
$ws.Cells.Item(6, 20).Value = 0.01543999301549823
$ws.Cells.Item(7, 7).Value = 0.178715
$ws.Cells.Item(7, 8).Value = 0.536145
$ws.Cells.Item(7, 9).Value = 0.09904930989061336
$ws.Cells.Item(7, 10).Value = 0.09904930989061336
$ws.Cells.Item(7, 13).Value = 19.66149466666667
$ws.Cells.Item(7, 14).Value = 58.984484
$ws.Cells.Item(7, 15).Value = 0.1388032745500265
$ws.Cells.Item(7, 16).Value = 0.1388032745500265
$ws.Cells.Item(7, 17).Value = 3.513804019353333
$ws.Cells.Item(7, 18).Value = 31.62423617418
$ws.Cells.Item(7, 19).Value = 0.01374836855473746
$ws.Cells.Item(7, 20).Value = 0.01374836855473746
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.450498333333333
$ws.Cells.Item(8, 8).Value = 4.351495
$ws.Cells.Item(8, 9).Value = 0.8039104659046613
$ws.Cells.Item(8, 10).Value = 0.8039104659046612
$ws.Cells.Item(8, 13).Value = 36.81180933333333
$ws.Cells.Item(8, 14).Value = 110.435428
$ws.Cells.Item(8, 15).Value = 0.2598784967371026
$ws.Cells.Item(8, 16).Value = 0.2598784967371026
$ws.Cells.Item(8, 17).Value = 53.39546808498444
$ws.Cells.Item(8, 18).Value = 480.5592127648599
$ws.Cells.Item(8, 19).Value = 0.2089190433905272
$ws.Cells.Item(8, 20).Value = 0.2089190433905271
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.450498333333333
$ws.Cells.Item(9, 8).Value = 4.351495
$ws.Cells.Item(9, 9).Value = 0.8039104659046613
$ws.Cells.Item(9, 10).Value = 0.8039104659046612
$ws.Cells.Item(9, 15).Value = 0.1970278712683331
$ws.Cells.Item(9, 16).Value = 0.197027871268333
$ws.Cells.Item(9, 17).Value = 40.48197732497778
$ws.Cells.Item(9, 18).Value = 364.3377959248
$ws.Cells.Item(9, 19).Value = 0.1583927677875293
$ws.Cells.Item(9, 20).Value = 0.1583927677875292
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.450498333333333
$ws.Cells.Item(10, 8).Value = 4.351495
$ws.Cells.Item(10, 9).Value = 0.8039104659046613
$ws.Cells.Item(10, 10).Value = 0.8039104659046612
$ws.Cells.Item(10, 13).Value = 21.95609833333333
$ws.Cells.Item(10, 14).Value = 65.868295
$ws.Cells.Item(10, 15).Value = 0.1550023737603119
$ws.Cells.Item(10, 16).Value = 0.1550023737603119
$ws.Cells.Item(10, 17).Value = 31.84728403900278
$ws.Cells.Item(10, 18).Value = 286.625556351025
$ws.Cells.Item(10, 19).Value = 0.1246080305059808
$ws.Cells.Item(10, 20).Value = 0.1246080305059808
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.450498333333333
$ws.Cells.Item(11, 8).Value = 4.351495
$ws.Cells.Item(11, 9).Value = 0.8039104659046613
$ws.Cells.Item(11, 10).Value = 0.8039104659046612
$ws.Cells.Item(11, 13).Value = 13.23098133333333
$ws.Cells.Item(11, 14).Value = 39.692944
$ws.Cells.Item(11, 15).Value = 0.09340609987756826
$ws.Cells.Item(11, 16).Value = 0.09340609987756825
$ws.Cells.Item(11, 17).Value = 19.19151637236444
$ws.Cells.Item(11, 18).Value = 172.72364735128
$ws.Cells.Item(11, 19).Value = 0.07509014127091322
$ws.Cells.Item(11, 20).Value = 0.07509014127091321
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.450498333333333
$ws.Cells.Item(12, 8).Value = 4.351495
$ws.Cells.Item(12, 9).Value = 0.8039104659046613
$ws.Cells.Item(12, 10).Value = 0.8039104659046612
$ws.Cells.Item(12, 13).Value = 22.080681
$ws.Cells.Item(12, 14).Value = 66.242043
$ws.Cells.Item(12, 15).Value = 0.1558818838066577
$ws.Cells.Item(12, 16).Value = 0.1558818838066577
$ws.Cells.Item(12, 17).Value = 32.02799098936499
$ws.Cells.Item(12, 18).Value = 288.251918904285
$ws.Cells.Item(12, 19).Value = 0.1253150778371065
$ws.Cells.Item(12, 20).Value = 0.1253150778371065
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.450498333333333
$ws.Cells.Item(13, 8).Value = 4.351495
$ws.Cells.Item(13, 9).Value = 0.8039104659046613
$ws.Cells.Item(13, 10).Value = 0.8039104659046612
$ws.Cells.Item(13, 13).Value = 19.66149466666667
$ws.Cells.Item(13, 14).Value = 58.984484
$ws.Cells.Item(13, 15).Value = 0.1388032745500265
$ws.Cells.Item(13, 16).Value = 0.1388032745500265
$ws.Cells.Item(13, 17).Value = 28.51896524484222
$ws.Cells.Item(13, 18).Value = 256.67068720358
$ws.Cells.Item(13, 19).Value = 0.1115854051126044
$ws.Cells.Item(13, 20).Value = 0.1115854051126044
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.17509
$ws.Cells.Item(14, 8).Value = 0.52527
$ws.Cells.Item(14, 9).Value = 0.09704022420472538
$ws.Cells.Item(14, 10).Value = 0.09704022420472537
$ws.Cells.Item(14, 13).Value = 36.81180933333333
$ws.Cells.Item(14, 14).Value = 110.435428
$ws.Cells.Item(14, 15).Value = 0.2598784967371026
$ws.Cells.Item(14, 16).Value = 0.2598784967371026
$ws.Cells.Item(14, 17).Value = 6.445379696173332
$ws.Cells.Item(14, 18).Value = 58.00841726556
$ws.Cells.Item(14, 19).Value = 0.02521866758935543
$ws.Cells.Item(14, 20).Value = 0.02521866758935543
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.17509
$ws.Cells.Item(15, 8).Value = 0.52527
$ws.Cells.Item(15, 9).Value = 0.09704022420472538
$ws.Cells.Item(15, 10).Value = 0.09704022420472537
$ws.Cells.Item(15, 15).Value = 0.1970278712683331
$ws.Cells.Item(15, 16).Value = 0.197027871268333
$ws.Cells.Item(15, 17).Value = 4.886589144533334
$ws.Cells.Item(15, 18).Value = 43.9793023008
$ws.Cells.Item(15, 19).Value = 0.01911962880245881
$ws.Cells.Item(15, 20).Value = 0.01911962880245881
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.17509
$ws.Cells.Item(16, 8).Value = 0.52527
$ws.Cells.Item(16, 9).Value = 0.09704022420472538
$ws.Cells.Item(16, 10).Value = 0.09704022420472537
$ws.Cells.Item(16, 13).Value = 21.95609833333333
$ws.Cells.Item(16, 14).Value = 65.868295
$ws.Cells.Item(16, 15).Value = 0.1550023737603119
$ws.Cells.Item(16, 16).Value = 0.1550023737603119
$ws.Cells.Item(16, 17).Value = 3.844293257183333
$ws.Cells.Item(16, 18).Value = 34.59863931465
$ws.Cells.Item(16, 19).Value = 0.01504146510196531
$ws.Cells.Item(16, 20).Value = 0.01504146510196531
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.17509
$ws.Cells.Item(17, 8).Value = 0.52527
$ws.Cells.Item(17, 9).Value = 0.09704022420472538
$ws.Cells.Item(17, 10).Value = 0.09704022420472537
$ws.Cells.Item(17, 13).Value = 13.23098133333333
$ws.Cells.Item(17, 14).Value = 39.692944
$ws.Cells.Item(17, 15).Value = 0.09340609987756826
$ws.Cells.Item(17, 16).Value = 0.09340609987756825
$ws.Cells.Item(17, 17).Value = 2.316612521653333
$ws.Cells.Item(17, 18).Value = 20.84951269488
$ws.Cells.Item(17, 19).Value = 0.009064148874208196
$ws.Cells.Item(17, 20).Value = 0.009064148874208193
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.17509
$ws.Cells.Item(18, 8).Value = 0.52527
$ws.Cells.Item(18, 9).Value = 0.09704022420472538
$ws.Cells.Item(18, 10).Value = 0.09704022420472537
$ws.Cells.Item(18, 13).Value = 22.080681
$ws.Cells.Item(18, 14).Value = 66.242043
$ws.Cells.Item(18, 15).Value = 0.1558818838066577
$ws.Cells.Item(18, 16).Value = 0.1558818838066577
$ws.Cells.Item(18, 17).Value = 3.86610643629
$ws.Cells.Item(18, 18).Value = 34.79495792661
$ws.Cells.Item(18, 19).Value = 0.01512681295405302
$ws.Cells.Item(18, 20).Value = 0.01512681295405301
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.17509
$ws.Cells.Item(19, 8).Value = 0.52527
$ws.Cells.Item(19, 9).Value = 0.09704022420472538
$ws.Cells.Item(19, 10).Value = 0.09704022420472537
$ws.Cells.Item(19, 13).Value = 19.66149466666667
$ws.Cells.Item(19, 14).Value = 58.984484
$ws.Cells.Item(19, 15).Value = 0.1388032745500265
$ws.Cells.Item(19, 16).Value = 0.1388032745500265
$ws.Cells.Item(19, 17).Value = 3.442531101186666
$ws.Cells.Item(19, 18).Value = 30.98277991068
$ws.Cells.Item(19, 19).Value = 0.1115854051126044
$ws.Cells.Item(19, 20).Value = 0.1115854051126044
